$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells for the columns that Excel would otherwise
# auto-detect as numbers/dates (Caso, F. De Reclamo, Comuna, OT), so the
# values stay literal text like every other row in this column.
$ws.Range("A84:B85").NumberFormat = "@"
$ws.Range("D84:E85").NumberFormat = "@"

# Row 84
$ws.Range("A84").Value = "-583"
$ws.Range("B84").Value = "9/8/2025"
$ws.Range("C84").Value = "Av Eva Perón 1145"
$ws.Range("D84").Value = "7"
$ws.Range("E84").Value = "809504290"
$ws.Range("F84").Value = "AYKO"
$ws.Range("G84").Value = "Pendiente"
$ws.Range("H84").Value = "Picada"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = "Cambio"
$ws.Range("K84").Value = "Sin equipos"
$ws.Range("L84").Value = "Pasante"
$ws.Range("M84").Value = -58.441547
$ws.Range("N84").Value = -34.630481
$ws.Range("O84").Value = "Boedo"
$ws.Range("P84").Value = "Capital Sur"

# Row 85
$ws.Range("A85").Value = "-584"
$ws.Range("B85").Value = "9/8/2025"
$ws.Range("C85").Value = "General Gregorio Aráoz de Lamadrid 865"
$ws.Range("D85").Value = "4"
$ws.Range("E85").Value = "809504300"
$ws.Range("F85").Value = "AYKO"
$ws.Range("G85").Value = "Pendiente"
$ws.Range("H85").Value = "Columna colgando"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = "Cambio"
$ws.Range("K85").Value = "Sin equipos"
$ws.Range("L85").Value = "Pasante"
$ws.Range("M85").Value = -58.364566
$ws.Range("N85").Value = -34.639404
$ws.Range("O85").Value = "San Telmo"
$ws.Range("P85").Value = "Capital Sur"
